$d = $word.ActiveDocument

# Heading paragraph indices (1-based) mapped to their TOC bookmark name.
$targets = @(
    @{ Index = 1;  Name = "_Toc5728364" },   # Story
    @{ Index = 2;  Name = "_Toc5728365" },   # Synopsis
    @{ Index = 4;  Name = "_Toc5728366" },   # Complete story
    @{ Index = 6;  Name = "_Toc5728367" },   # Backstory
    @{ Index = 8;  Name = "_Toc5728368" },   # Narrative devices
    @{ Index = 10; Name = "_Toc5728369" }    # Subplots
)

foreach ($t in $targets) {
    $p = $d.Paragraphs($t.Index)
    $r = $p.Range
    # Exclude the trailing paragraph mark so the bookmark wraps only the
    # heading text (matches bookmarkStart/bookmarkEnd placed around the run).
    $r.End = $r.End - 1
    $d.Bookmarks.Add($t.Name, $r)
}

Write-Output "done"
